$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue "I3" "1"
Set-TextValue "I4" "1"
Set-TextValue "I5" "2"
Set-TextValue "I6" "2"
Set-TextValue "I10" "1"
Set-TextValue "I11" "2"
Set-TextValue "I15" "1"
